$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Based on threshold review, the marker for these samples changed from G418 to NAT
$ws.Range("J20").Value = "NAT"
$ws.Range("J21").Value = "NAT"
$ws.Range("J22").Value = "NAT"
